$d = $word.ActiveDocument

# Locate the paragraph that starts the footer block to be removed
# ("Ver no Jupiter Salvar em pdf Salvar em docx").
$findStart = $d.Content
$startFound = $findStart.Find.Execute(
    "Ver no Jupiter Salvar em pdf Salvar em docx",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startPara = $findStart.Paragraphs(1)

# Locate the copyright paragraph
# ("© 2020 . Contact: ... Powered by Jekyll and Github pages. ...").
$findEnd = $d.Content
$endFound = $findEnd.Find.Execute(
    "Powered by Jekyll and Github pages",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$copyrightPara = $findEnd.Paragraphs(1)

# The blank paragraph immediately following the copyright line is also
# removed, leaving only the blank paragraph that already separated the
# bibliography entry from this footer block.
$endPara = $copyrightPara.Next()

$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()
